$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.417.82"
$ws.Range("E2").Value = "  +3.32%  "

$ws.Range("D3").Value = "1.796.30"
$ws.Range("E3").Value = "  +4.04%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.04"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3814"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3458"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.75%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.25"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.207"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07545"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.16"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +10.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.507"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.42%  "

$ws.Range("D15").Value = "1.795.04"
$ws.Range("E15").Value = "  +4.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.104"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001103"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06672"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.20"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.537"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.52%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.42"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.18%  "

$ws.Range("D23").Value = "27.383.14"
$ws.Range("E23").Value = "  +3.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.59"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.436"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.596"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.507"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.48"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +10.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "152.32"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.08%  "

$ws.Range("D30").Value = "1.997.21"
$ws.Range("E30").Value = "  +4.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.18"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.070"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.156"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.74%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08731"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.34"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.691"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.480"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6939"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +11.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.922"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06389"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2211"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02343"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.278"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.41"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.90%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6496"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.876"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.133"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.67"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07192"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.47"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.09%  "
